# edit.ps1
# Applies the "Doing Updates for Financials" change:
#  - Inserts two new columns before column D (new quarterly data for
#    quarters ending 2018-12-31 and 2018-09-30)
#  - Shifts the previously-existing quarterly data two columns to the right
#    (this happens automatically via the column Insert)
#  - Fills in the values for the two newly inserted columns
#  - Applies a handful of historical restatements that occurred in the
#    same update (values that changed in the now-shifted older columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank columns before column D. This shifts all existing
#    D:K data to F:M automatically.
$ws.Range("D:E").Insert()

# 2. The newly inserted D:E columns come out formatted like column C
#    (format is copied from the column to the left by default). Column F
#    (the original column D) has the correct per-row number/date format,
#    so copy formats from F into D:E to match. This is done per
#    contiguous block of rows that actually contain data in D:M (skipping
#    the section-header rows 5, 6, 37 and 79, which have no D:M cells at
#    all, either before or after the edit).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$newData = @(
    @{ Row=7; D=43465; E=43373 },
    @{ Row=8; D=12400; E=11400 },
    @{ Row=9; D=1600; E=1500 },
    @{ Row=10; D=10800; E=9900 },
    @{ Row=11; D=$null; E=$null },
    @{ Row=12; D="NA"; E="NA" },
    @{ Row=13; D=0; E=0 },
    @{ Row=14; D=0; E=0 },
    @{ Row=15; D=0; E=0 },
    @{ Row=16; D=$null; E=$null },
    @{ Row=17; D=23600; E=21800 },
    @{ Row=18; D=-11200; E=-10400 },
    @{ Row=19; D=$null; E=$null },
    @{ Row=20; D=0; E=100 },
    @{ Row=21; D=-9300; E=-8500 },
    @{ Row=22; D=0; E=0 },
    @{ Row=23; D=-11200; E=-10300 },
    @{ Row=24; D=-100; E=0 },
    @{ Row=25; D=0; E=0 },
    @{ Row=26; D=-11100; E=-10300 },
    @{ Row=27; D=-11100; E=-10300 },
    @{ Row=28; D=0; E=0 },
    @{ Row=29; D=0; E="NA" },
    @{ Row=30; D=0; E=0 },
    @{ Row=31; D=0; E=0 },
    @{ Row=32; D=0; E=-100 },
    @{ Row=33; D=-11100; E=-10300 },
    @{ Row=34; D=0; E=0 },
    @{ Row=35; D=-11100; E=-10300 },
    @{ Row=38; D=43465; E=43373 },
    @{ Row=39; D=$null; E=$null },
    @{ Row=40; D=$null; E=$null },
    @{ Row=41; D=30000; E=30800 },
    @{ Row=42; D=0; E=0 },
    @{ Row=43; D=1300; E=1500 },
    @{ Row=44; D=0; E=0 },
    @{ Row=45; D=3200; E=3200 },
    @{ Row=46; D=34500; E=35500 },
    @{ Row=47; D=12700; E=12700 },
    @{ Row=48; D=21700; E=20300 },
    @{ Row=49; D=38200; E=36400 },
    @{ Row=50; D=0; E=0 },
    @{ Row=51; D=0; E=0 },
    @{ Row=52; D=0; E=0 },
    @{ Row=53; D=0; E=0 },
    @{ Row=54; D=107200; E=105000 },
    @{ Row=55; D=$null; E=$null },
    @{ Row=56; D=$null; E=$null },
    @{ Row=57; D=8000; E=8200 },
    @{ Row=58; D=0; E=0 },
    @{ Row=59; D=5000; E=4500 },
    @{ Row=60; D=13000; E=12800 },
    @{ Row=61; D=12500; E=0 },
    @{ Row=62; D=0; E=200 },
    @{ Row=63; D=0; E=0 },
    @{ Row=64; D=0; E=0 },
    @{ Row=65; D=0; E=0 },
    @{ Row=66; D=25700; E=12900 },
    @{ Row=67; D=$null; E=$null },
    @{ Row=68; D=0; E=0 },
    @{ Row=69; D=0; E=0 },
    @{ Row=70; D=0; E=0 },
    @{ Row=71; D=0; E=0 },
    @{ Row=72; D=-58200; E=-47100 },
    @{ Row=73; D=0; E=0 },
    @{ Row=74; D=0; E=0 },
    @{ Row=75; D=0; E=0 },
    @{ Row=76; D=81500; E=92000 },
    @{ Row=77; D=0; E=0 },
    @{ Row=80; D=43465; E=43373 },
    @{ Row=81; D=-11100; E=-10300 },
    @{ Row=82; D=$null; E=$null },
    @{ Row=83; D=1900; E=1800 },
    @{ Row=84; D=0; E=0 },
    @{ Row=85; D=0; E=0 },
    @{ Row=86; D=0; E=0 },
    @{ Row=87; D=0; E=0 },
    @{ Row=88; D=0; E=0 },
    @{ Row=89; D=-8200; E=-6300 },
    @{ Row=90; D=$null; E=$null },
    @{ Row=91; D=-5100; E=-4300 },
    @{ Row=92; D=0; E=0 },
    @{ Row=93; D=0; E=0 },
    @{ Row=94; D=-5100; E=-4300 },
    @{ Row=95; D=$null; E=$null },
    @{ Row=96; D=0; E=0 },
    @{ Row=97; D=0; E=0 },
    @{ Row=98; D=0; E=0 },
    @{ Row=99; D=0; E=0 },
    @{ Row=100; D=12500; E=300 },
    @{ Row=101; D=0; E=0 },
    @{ Row=102; D=-800; E=-10300 },
)


# 3. Write the new quarterly values into columns D and E for every row.
foreach ($item in $newData) {
    $r = $item.Row
    if ($null -ne $item.D) {
        $ws.Cells.Item($r, 4).Value = $item.D
    }
    if ($null -ne $item.E) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}

$restatedH = @(
    @{ Row=48; Value=34100 },
    @{ Row=57; Value=4300 },
    @{ Row=60; Value=7700 },
    @{ Row=61; Value=12500 },
    @{ Row=91; Value=-3000 },
)

# 4. A handful of values in the historical data (now in column H, which
#    was column F before the insert) were restated as part of this update.
foreach ($item in $restatedH) {
    $ws.Cells.Item($item.Row, 8).Value = $item.Value
}
